# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Asura_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 5715.353
$ws.Range("I98").Value = 4470.1787
$ws.Range("J98").Value = 11526.167
$ws.Range("K98").Value = 4470.1787
$ws.Range("L98").Value = 11526.167
$ws.Range("M98").Value = -2972.1787
$ws.Range("N98").Value = -14522.167
$ws.Range("H122").Value = 5715.353
$ws.Range("I122").Value = 4470.1787
$ws.Range("J122").Value = 11526.167
$ws.Range("K122").Value = 13410.5361
$ws.Range("L122").Value = 34578.501
$ws.Range("M122").Value = -10960.5361
$ws.Range("N122").Value = -39478.501
$ws.Range("H132").Value = 1360.4579
$ws.Range("I132").Value = 1274.6836
$ws.Range("J132").Value = 3054.5
$ws.Range("K132").Value = 3824.0508
$ws.Range("L132").Value = 9163.5
$ws.Range("M132").Value = -1294.0508
$ws.Range("N132").Value = -14223.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280
$ws.Range("H140").Value = 146775
$ws.Range("J140").Value = 200000
$ws.Range("L140").Value = 200000
$ws.Range("N140").Value = -210360

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2314.2666
$ws.Range("I122").Value = 2293.6667
$ws.Range("J122").Value = 2396.6667
$ws.Range("K122").Value = 6881.000100000001
$ws.Range("L122").Value = 7190.000100000001
$ws.Range("M122").Value = -4431.000100000001
$ws.Range("N122").Value = -12090.0001
$ws.Range("H132").Value = 365304.06
$ws.Range("I132").Value = 401360.5
$ws.Range("K132").Value = 1204081.5
$ws.Range("M132").Value = -1201551.5
$ws.Range("H133").Value = 30095.445
$ws.Range("J133").Value = 30095.445
$ws.Range("L133").Value = 30095.445
$ws.Range("N133").Value = -35155.445

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 84590
$ws.Range("J132").Value = 84590
$ws.Range("L132").Value = 84590
$ws.Range("N132").Value = -94710
$ws.Range("H134").Value = 590703.9399999999
$ws.Range("I134").Value = 647101.0600000001
$ws.Range("J134").Value = 7933.3335
$ws.Range("K134").Value = 1941303.18
$ws.Range("L134").Value = 23800.0005
$ws.Range("M134").Value = -1938768.18
$ws.Range("N134").Value = -28870.0005

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1786.8334
$ws.Range("I122").Value = 1517.75
$ws.Range("J122").Value = 2325
$ws.Range("K122").Value = 4553.25
$ws.Range("L122").Value = 6975
$ws.Range("M122").Value = -2103.25
$ws.Range("N122").Value = -11875
$ws.Range("H132").Value = 1421
$ws.Range("I132").Value = 1421
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4263
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1733
$ws.Range("N132").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 770.1429000000001
$ws.Range("I8").Value = 770.1429000000001
$ws.Range("K8").Value = 2310.4287
$ws.Range("M8").Value = -2171.4287
$ws.Range("H131").Value = 824.2
$ws.Range("I131").Value = 420
$ws.Range("J131").Value = 859.34784
$ws.Range("K131").Value = 1260
$ws.Range("L131").Value = 2578.04352
$ws.Range("M131").Value = 3780
$ws.Range("N131").Value = -12658.04352
$ws.Range("H132").Value = 2212.611
$ws.Range("J132").Value = 2697.3
$ws.Range("L132").Value = 24275.7
$ws.Range("N132").Value = -29335.7
$ws.Range("H136").Value = 4041.5881
$ws.Range("J136").Value = 4389.8
$ws.Range("L136").Value = 13169.4
$ws.Range("N136").Value = -23369.4
$ws.Range("H139").Value = 2269.7896
$ws.Range("I139").Value = 1838.3334
$ws.Range("J139").Value = 3009.4285
$ws.Range("K139").Value = 5515.0002
$ws.Range("L139").Value = 9028.2855
$ws.Range("M139").Value = -375.0002000000004
$ws.Range("N139").Value = -19308.2855
$ws.Range("H141").Value = 2364.1428
$ws.Range("I141").Value = 2364.1428
$ws.Range("K141").Value = 7092.428400000001
$ws.Range("M141").Value = -1912.428400000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 28400
$ws.Range("J46").Value = 28400
$ws.Range("L46").Value = 28400
$ws.Range("N46").Value = -28712
$ws.Range("H102").Value = 2392
$ws.Range("I102").Value = 2515.3333
$ws.Range("J102").Value = 2207
$ws.Range("K102").Value = 2515.3333
$ws.Range("L102").Value = 2207
$ws.Range("M102").Value = -893.3332999999998
$ws.Range("N102").Value = -5451
$ws.Range("H122").Value = 3867.1162
$ws.Range("I122").Value = 3629.1353
$ws.Range("J122").Value = 5334.6665
$ws.Range("K122").Value = 10887.4059
$ws.Range("L122").Value = 16003.9995
$ws.Range("M122").Value = -8437.4059
$ws.Range("N122").Value = -20903.9995
$ws.Range("H132").Value = 1997.3143
$ws.Range("I132").Value = 1182.32
$ws.Range("K132").Value = 3546.96
$ws.Range("M132").Value = -1016.96

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4350.3125
$ws.Range("I7").Value = 4077.7778
$ws.Range("J7").Value = 4700.7144
$ws.Range("K7").Value = 4077.7778
$ws.Range("L7").Value = 4700.7144
$ws.Range("M7").Value = -3965.7778
$ws.Range("N7").Value = -4924.7144
$ws.Range("H40").Value = 9004.556
$ws.Range("I40").Value = 8171.091
$ws.Range("J40").Value = 10314.286
$ws.Range("K40").Value = 8171.091
$ws.Range("L40").Value = 10314.286
$ws.Range("M40").Value = -8035.091
$ws.Range("N40").Value = -10586.286
$ws.Range("H122").Value = 1019395.4
$ws.Range("J122").Value = 4967.1113
$ws.Range("L122").Value = 14901.3339
$ws.Range("N122").Value = -19801.3339
$ws.Range("H126").Value = 4350.3125
$ws.Range("I126").Value = 4077.7778
$ws.Range("J126").Value = 4700.7144
$ws.Range("K126").Value = 12233.3334
$ws.Range("L126").Value = 14102.1432
$ws.Range("M126").Value = -9763.3334
$ws.Range("N126").Value = -19042.1432

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 15000
$ws.Range("J48").Value = 15000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -16138
$ws.Range("H49").Value = 15000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H122").Value = 727.55554
$ws.Range("I122").Value = 528.2857
$ws.Range("J122").Value = 1425
$ws.Range("K122").Value = 1584.8571
$ws.Range("L122").Value = 4275
$ws.Range("M122").Value = 865.1428999999998
$ws.Range("N122").Value = -9175
$ws.Range("H126").Value = 7168.6924
$ws.Range("I126").Value = 8819.799999999999
$ws.Range("K126").Value = 26459.4
$ws.Range("M126").Value = -23989.4
$ws.Range("H136").Value = 1453.7693
$ws.Range("I136").Value = 1523.9524
$ws.Range("J136").Value = 1159
$ws.Range("K136").Value = 4571.857199999999
$ws.Range("L136").Value = 3477
$ws.Range("M136").Value = -2021.857199999999
$ws.Range("N136").Value = -8577

